$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1273.0159
$ws.Range("J112").Value = 1273.0159
$ws.Range("L112").Value = 3819.0477
$ws.Range("N112").Value = -6035.0477
$ws.Range("H116").Value = 5822.815
$ws.Range("I116").Value = 2450.3333
$ws.Range("J116").Value = 8520.799999999999
$ws.Range("K116").Value = 2450.3333
$ws.Range("L116").Value = 8520.799999999999
$ws.Range("M116").Value = 991.6667000000002
$ws.Range("N116").Value = -15404.8
$ws.Range("H125").Value = 1213
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1213
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 10917
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -15837
$ws.Range("H132").Value = 33802036
$ws.Range("I132").Value = 38617070
$ws.Range("K132").Value = 115851210
$ws.Range("M132").Value = -115848680
$ws.Range("H138").Value = 4760.441
$ws.Range("I138").Value = 2608.375
$ws.Range("J138").Value = 5098.0195
$ws.Range("K138").Value = 7825.125
$ws.Range("L138").Value = 15294.0585
$ws.Range("M138").Value = -2685.125
$ws.Range("N138").Value = -25574.0585
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1598.1818
$ws.Range("I61").Value = 1594
$ws.Range("J61").Value = 1609.3334
$ws.Range("K61").Value = 1594
$ws.Range("L61").Value = 1609.3334
$ws.Range("M61").Value = -1382
$ws.Range("N61").Value = -2033.3334
$ws.Range("H63").Value = 8660613
$ws.Range("I63").Value = 15392256
$ws.Range("J63").Value = 5642.857
$ws.Range("K63").Value = 15392256
$ws.Range("L63").Value = 5642.857
$ws.Range("M63").Value = -15391570
$ws.Range("N63").Value = -7014.857
$ws.Range("H66").Value = 8660613
$ws.Range("I66").Value = 15392256
$ws.Range("J66").Value = 5642.857
$ws.Range("K66").Value = 76961280
$ws.Range("L66").Value = 28214.285
$ws.Range("M66").Value = -76957848
$ws.Range("N66").Value = -35078.285
$ws.Range("H74").Value = 1505.8611
$ws.Range("I74").Value = 1007
$ws.Range("J74").Value = 6993.3335
$ws.Range("K74").Value = 1007
$ws.Range("L74").Value = 6993.3335
$ws.Range("M74").Value = -133
$ws.Range("N74").Value = -8741.333500000001
$ws.Range("H77").Value = 1505.8611
$ws.Range("I77").Value = 1007
$ws.Range("J77").Value = 6993.3335
$ws.Range("K77").Value = 5035
$ws.Range("L77").Value = 34966.6675
$ws.Range("M77").Value = -667
$ws.Range("N77").Value = -43702.6675
$ws.Range("H123").Value = 49643.2
$ws.Range("J123").Value = 49643.2
$ws.Range("L123").Value = 49643.2
$ws.Range("N123").Value = -59443.2
$ws.Range("H132").Value = 2338
$ws.Range("I132").Value = 1512.9524
$ws.Range("J132").Value = 4070.6
$ws.Range("K132").Value = 4538.857199999999
$ws.Range("L132").Value = 12211.8
$ws.Range("M132").Value = -2008.857199999999
$ws.Range("N132").Value = -17271.8
$ws.Range("H136").Value = 1598.1818
$ws.Range("I136").Value = 1594
$ws.Range("J136").Value = 1609.3334
$ws.Range("K136").Value = 4782
$ws.Range("L136").Value = 4828.0002
$ws.Range("M136").Value = -2232
$ws.Range("N136").Value = -9928.0002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 3859.8
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 3859.8
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 3859.8
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -4085.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8549069
$ws.Range("I16").Value = 18520182
$ws.Range("J16").Value = 2401.4285
$ws.Range("K16").Value = 18520182
$ws.Range("L16").Value = 2401.4285
$ws.Range("M16").Value = -18519895
$ws.Range("N16").Value = -2975.4285
$ws.Range("H31").Value = 6511.0376
$ws.Range("I31").Value = 2373.2778
$ws.Range("J31").Value = 8639.028
$ws.Range("K31").Value = 2373.2778
$ws.Range("L31").Value = 8639.028
$ws.Range("M31").Value = -2078.2778
$ws.Range("N31").Value = -9229.028
$ws.Range("H34").Value = 6511.0376
$ws.Range("I34").Value = 2373.2778
$ws.Range("J34").Value = 8639.028
$ws.Range("K34").Value = 2373.2778
$ws.Range("L34").Value = 8639.028
$ws.Range("M34").Value = -2171.2778
$ws.Range("N34").Value = -9043.028
$ws.Range("H113").Value = 8549069
$ws.Range("I113").Value = 18520182
$ws.Range("J113").Value = 2401.4285
$ws.Range("K113").Value = 18520182
$ws.Range("L113").Value = 2401.4285
$ws.Range("M113").Value = -18518012
$ws.Range("N113").Value = -6741.4285
$ws.Range("H134").Value = 4000.2195
$ws.Range("I134").Value = 4829.148
$ws.Range("J134").Value = 2401.5715
$ws.Range("K134").Value = 14487.444
$ws.Range("L134").Value = 7204.7145
$ws.Range("M134").Value = -11952.444
$ws.Range("N134").Value = -12274.7145
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 2746.4443
$ws.Range("I49").Value = 2529.4285
$ws.Range("J49").Value = 3506
$ws.Range("K49").Value = 7588.2855
$ws.Range("L49").Value = 10518
$ws.Range("M49").Value = -7432.2855
$ws.Range("N49").Value = -10830
$ws.Range("H64").Value = 999.5
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 999.5
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H75").Value = 3507
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 3507
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 10521
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -12517
$ws.Range("H78").Value = 3507
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 3507
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 31563
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -41547
$ws.Range("H103").Value = 2505
$ws.Range("I103").Value = 525
$ws.Range("K103").Value = 1575
$ws.Range("M103").Value = -696
$ws.Range("H114").Value = 1600
$ws.Range("I114").Value = 675
$ws.Range("J114").Value = 2833.3333
$ws.Range("K114").Value = 2025
$ws.Range("L114").Value = 8499.999899999999
$ws.Range("M114").Value = 1229
$ws.Range("N114").Value = -15007.9999
$ws.Range("H117").Value = 1387
$ws.Range("I117").Value = 329
$ws.Range("J117").Value = 1916
$ws.Range("K117").Value = 987
$ws.Range("L117").Value = 5748
$ws.Range("M117").Value = 2455
$ws.Range("N117").Value = -12632
$ws.Range("H121").Value = 1901.4333
$ws.Range("I121").Value = 348
$ws.Range("J121").Value = 2042.6545
$ws.Range("K121").Value = 1044
$ws.Range("L121").Value = 6127.9635
$ws.Range("M121").Value = 266
$ws.Range("N121").Value = -8747.9635
$ws.Range("H129").Value = 3162.1052
$ws.Range("I129").Value = 2531.111
$ws.Range("J129").Value = 3730
$ws.Range("K129").Value = 7593.333
$ws.Range("L129").Value = 11190
$ws.Range("M129").Value = -2593.333
$ws.Range("N129").Value = -21190
$ws.Range("H131").Value = 798.54736
$ws.Range("J131").Value = 825.13336
$ws.Range("L131").Value = 2475.40008
$ws.Range("N131").Value = -12555.40008
$ws.Range("H132").Value = 2751.4736
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 7942.143
$ws.Range("I2").Value = 863
$ws.Range("J2").Value = 13251.5
$ws.Range("K2").Value = 863
$ws.Range("L2").Value = 13251.5
$ws.Range("M2").Value = -751
$ws.Range("N2").Value = -13475.5
$ws.Range("H7").Value = 3988.5789
$ws.Range("I7").Value = 1518.3
$ws.Range("K7").Value = 1518.3
$ws.Range("M7").Value = -1406.3
$ws.Range("H74").Value = 47960
$ws.Range("J74").Value = 47960
$ws.Range("L74").Value = 47960
$ws.Range("N74").Value = -49956
$ws.Range("H77").Value = 47960
$ws.Range("J77").Value = 47960
$ws.Range("L77").Value = 143880
$ws.Range("N77").Value = -153864
$ws.Range("H103").Value = 34900.5
$ws.Range("J103").Value = 34900.5
$ws.Range("L103").Value = 34900.5
$ws.Range("N103").Value = -37244.5
$ws.Range("H110").Value = 27883.857
$ws.Range("J110").Value = 27883.857
$ws.Range("L110").Value = 27883.857
$ws.Range("N110").Value = -36063.857
$ws.Range("H126").Value = 3988.5789
$ws.Range("I126").Value = 1518.3
$ws.Range("K126").Value = 4554.9
$ws.Range("M126").Value = -2084.9
$ws.Range("H136").Value = 6252.933
$ws.Range("I136").Value = 1915.6666
$ws.Range("K136").Value = 5746.9998
$ws.Range("M136").Value = -3196.9998
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 3040
$ws.Range("I3").Value = 2400
$ws.Range("K3").Value = 2400
$ws.Range("M3").Value = -2286
$ws.Range("H9").Value = 6500
$ws.Range("J9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("N9").Value = -10280
$ws.Range("H119").Value = 26349
$ws.Range("J119").Value = 26349
$ws.Range("L119").Value = 26349
$ws.Range("N119").Value = -36025
$ws.Range("H125").Value = 43772
$ws.Range("J125").Value = 43772
$ws.Range("L125").Value = 43772
$ws.Range("N125").Value = -53612
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 15876232
$ws.Range("I132").Value = 2563
$ws.Range("J132").Value = 30306840
$ws.Range("K132").Value = 7689
$ws.Range("L132").Value = 90920520
$ws.Range("M132").Value = -5159
$ws.Range("N132").Value = -90925580
